# StyleTagTemplate.xlsx edit:
#  - add a new worksheet "class only" (sheetId 14) after the last existing
#    sheet ("width-height"), containing 5 new jt:style "class" examples
#  - make the new sheet the active/selected tab (moves tabSelected from the
#    first sheet to this new one, and updates the workbook's activeTab)

$wb = $excel.ActiveWorkbook

# --- create the new sheet at the end of the tab strip -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "class only"

# --- column B width (chars), to roughly match the authored template ---
# (36.5703125 in stored OOXML units; the ColumnWidth setter here quantizes
# to 1/6-character steps, so 36.5 is the closest reachable stored value)
$ws.Columns.Item(2).ColumnWidth = 35.67

# --- cell contents ------------------------------------------------------
# Shared-string insertion order matters (it determines the index each
# string gets in xl/sharedStrings.xml), so write them in the same order
# the template does: redBoxCenter, doesNotExist, blueBoldBigText,
# blueBoldBigText-with-overrides, then the combined-classes example.
$ws.Range("B2").Value = '<jt:style class="redBoxCenter">redBoxCenter</jt:style>'
$ws.Range("B6").Value = '<jt:style class="doesNotExist">doesNotExist</jt:style>'
$ws.Range("B4").Value = '<jt:style class="blueBoldBigText">blueBoldBigText</jt:style>'
$ws.Range("B8").Value = '<jt:style class="blueBoldBigText" style="font-color: green; font-weight: normal; font-italic: true">blueBoldBigText with green non-bold italic</jt:style>'
$ws.Range("B10").Value = '<jt:style class="redBoxCenter; blueBoldBigText">redBoxCenter; blueBoldBigText</jt:style>'

# --- make the new sheet the active tab ----------------------------------
# This clears tabSelected on whichever sheet had it before ("alignment")
# and sets it (plus the workbook's activeTab index) on the new sheet.
$ws.Activate()
